$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AC3").Value = 11
$ws.Range("AP3").Value = 19
$ws.Range("AU3").Value = 9

$ws.Range("S4").Value = 1.53
$ws.Range("T4").Value = 2.38

$ws.Range("S5").Value = 1.57
